$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# Update the "version" column (E) for Disease Ontology (row 3) and
# Experimental Factor Ontology (row 4) to reflect the Jan 2023 DO/EFO update.
$ws.Range("E3").Value = "v2023-01-31"
$ws.Range("E4").Value = "v3.50.0"

# Update the active selection to reflect the last-edited cell.
$ws.Range("E5").Select()
